# Apply updated crypto price/volume data to Sheet1 (GitHub Actions scrape refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "306.41") are temporarily forced to Text format so the written value
# stays a string, matching the source data - then the style is reset to
# "Normal" so no stray number-format survives on the cell.
$textCells = @("D5", "D6", "D9", "D10", "D11", "D14", "D15", "D17", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D33", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '39.652.84'
$ws.Range("E2").Value = '  -2.96%  '
$ws.Range("D3").Value = '2.308.60'
$ws.Range("E3").Value = '  -4.22%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '306.41'
$ws.Range("E5").Value = '  -2.94%  '
$ws.Range("D6").Value = '81.60'
$ws.Range("E6").Value = '  -7.57%  '
$ws.Range("E7").Value = '  -3.60%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '0.470'
$ws.Range("E9").Value = '  -4.70%  '
$ws.Range("D10").Value = '0.0783'
$ws.Range("E10").Value = '  -5.43%  '
$ws.Range("D11").Value = '28.65'
$ws.Range("E11").Value = '  -9.35%  '
$ws.Range("E12").Value = '  -0.22%  '
$ws.Range("D13").Value = '2.670.77'
$ws.Range("E13").Value = '  -4.13%  '
$ws.Range("D14").Value = '6.21'
$ws.Range("E14").Value = '  -6.60%  '
$ws.Range("D15").Value = '14.39'
$ws.Range("E15").Value = '  -7.47%  '
$ws.Range("D16").Value = '2.329.75'
$ws.Range("E16").Value = '  -4.19%  '
$ws.Range("D17").Value = '0.736'
$ws.Range("E17").Value = '  -4.30%  '
$ws.Range("D18").Value = '39.586.10'
$ws.Range("E18").Value = '  -2.94%  '
$ws.Range("D19").Value = '0.0₃0880'
$ws.Range("E19").Value = '  -4.11%  '
$ws.Range("E20").Value = '  -4.85%  '
$ws.Range("D21").Value = '67.31'
$ws.Range("E21").Value = '  -5.94%  '
$ws.Range("D22").Value = '10.22'
$ws.Range("E22").Value = '  -5.55%  '
$ws.Range("D23").Value = '231.72'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = '1.00'
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("B25").Value = 'PancakeSwap'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D25").Value = '2.47'
$ws.Range("E25").Value = '  -7.06%  '
$ws.Range("D26").Value = '1.77'
$ws.Range("E26").Value = '  -4.26%  '
$ws.Range("D27").Value = '22.76'
$ws.Range("E27").Value = '  -4.85%  '
$ws.Range("E28").Value = '  -1.93%  '
$ws.Range("D29").Value = '9.02'
$ws.Range("E29").Value = '  -4.90%  '
$ws.Range("D30").Value = '32.41'
$ws.Range("E30").Value = '  -5.13%  '
$ws.Range("D31").Value = '150.52'
$ws.Range("E31").Value = '  -4.11%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("D33").Value = '4.94'
$ws.Range("E33").Value = '  -5.17%  '
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").Value = '0.0699'
$ws.Range("E35").Value = '  -4.99%  '
$ws.Range("E36").Value = '  -1.39%  '
$ws.Range("D37").Value = '2.69'
$ws.Range("E37").Value = '  -6.50%  '
$ws.Range("D38").Value = '0.0960'
$ws.Range("E38").Value = '  -3.70%  '
$ws.Range("D39").Value = '15.10'
$ws.Range("E39").Value = '  -7.85%  '
$ws.Range("D40").Value = '1.63'
$ws.Range("E40").Value = '  -7.49%  '
$ws.Range("D41").Value = '3.66'
$ws.Range("E41").Value = '  -4.67%  '
$ws.Range("D42").Value = '2.27'
$ws.Range("E42").Value = '  -1.05%  '
$ws.Range("D43").Value = '1.958.02'
$ws.Range("E43").Value = '  -1.34%  '
$ws.Range("D44").Value = '0.0256'
$ws.Range("E44").Value = '  -5.83%  '
$ws.Range("D45").Value = '16.64'
$ws.Range("E45").Value = '  -9.65%  '
$ws.Range("D46").Value = '9.19'
$ws.Range("E46").Value = '  -1.35%  '
$ws.Range("D47").Value = '2.59'
$ws.Range("E47").Value = '  -8.91%  '
$ws.Range("D48").Value = '2.532.62'
$ws.Range("E48").Value = '  -4.29%  '
$ws.Range("D49").Value = '90.46'
$ws.Range("E49").Value = '  -3.50%  '
$ws.Range("D50").Value = '68.34'
$ws.Range("E50").Value = '  -6.24%  '
$ws.Range("D51").Value = '48.26'
$ws.Range("E51").Value = '  -5.89%  '

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
